# Apply the "usb-copy" preset row edit to the presets worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new preset row (row 4) ---
$ws.Range("A4").Value = "Comosite #1"
$ws.Range("B4").Value = "50-90"
$ws.Range("C4").Value = "50,70,90"
$ws.Range("D4").Value = 60120120
$ws.Range("D4").NumberFormat = "#,##0"

# --- Nudge the cached outline-level high-water marks that Excel keeps in
#     sheetFormatPr (outlineLevelRow/outlineLevelCol) up to 3, matching the
#     target file, without leaving a visible/used row or column grouped. ---
$ws.Columns.Item(5).OutlineLevel = 3
$ws.Rows.Item(6).OutlineLevel = 3
$ws.Rows.Item(6).Delete()

# --- Update the selection to match the post-edit state ---
$ws.Range("D5").Select()
